$d = $word.ActiveDocument

# --------------------------------------------------------------------
# 1. Insert a new title paragraph ("{{ title }}") at the very start of
#    the document body, before the existing "{{ image }}" paragraph.
#    InsertParagraphBefore() clones the formatting (pPr/rPr) of the
#    paragraph it is attached to, which is exactly the styling the
#    target markup uses (center-justified, sz 44 run props).
# --------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertParagraphBefore()
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.Text = "{{ title }}"

# --------------------------------------------------------------------
# 2. Collapse the run-per-token markup in the two data-row table cells
#    into a single run containing Jinja-style conditional text.
#    Each cell currently holds 9 separate <w:r> runs; we rewrite the
#    cell's paragraph so it ends up with exactly one run (keeping an
#    empty <w:rPr/> on it, matching the target markup).
# --------------------------------------------------------------------
function Get-ParaIndexContaining($searchText) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($searchText)) {
            return $i
        }
    }
    return -1
}

function Set-CellTemplateText($paraIndex, $lastRunText, $newText) {
    $para = $d.Paragraphs.Item($paraIndex)
    $full = $para.Range
    $full.MoveEnd(1, -1)   # exclude the paragraph mark
    $paraStart = $full.Start
    $paraEnd = $full.End

    # Delete everything except the final run (it carries no
    # xml:space="preserve", unlike the earlier runs which do - keeping
    # it avoids picking up a stray xml:space attribute on the rebuilt
    # run).
    $lastRunLen = $lastRunText.Length
    if (($paraEnd - $paraStart) -gt $lastRunLen) {
        $pre = $d.Range($paraStart, $paraEnd - $lastRunLen)
        $pre.Delete()
    }

    $remaining = $d.Paragraphs.Item($paraIndex).Range
    $remaining.MoveEnd(1, -1)
    $remaining.Text = $newText
}

$idxRow0 = Get-ParaIndexContaining "row[0][0]"
Set-CellTemplateText $idxRow0 "row[0][2] }}" "{% if row[0] %}{{ row[0][0] }} {% if number %} - {{ row[0][1] }}{% endif %}{% if percent %} - %{{ row[0][2] }}{% endif %}{% endif %}"

$idxRow1 = Get-ParaIndexContaining "row[1][0]"
Set-CellTemplateText $idxRow1 "row[1][2] }}" "{% if row[1] %}{{ row[1][0] }} {% if number %} - {{ row[1][1] }}{% endif %}{% if percent %} - %{{ row[1][2] }}{% endif %}{% endif %}"
